$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.077887892723083
$ws.Range("B1").Value = 1.506378412246704
$ws.Range("C1").Value = 3.906706094741821
$ws.Range("D1").Value = 1.7109534740448
$ws.Range("E1").Value = 0.9360960125923157
